# Auto-generated Excel COM-interop script
# Applies value updates across sheets "展览", "演出", "本地生活", "全部类型"
# matching the described diff (refreshed scrape snapshot values).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 1151
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202409/rDtevgY01726641860615.jpeg"
$ws.Range("F9").Value = 230
$ws.Range("F11").Value = 8027
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = 9540
$ws.Range("F17").Value = 468
$ws.Range("F25").Value = 46
$ws.Range("F28").Value = 1638
$ws.Range("F30").Value = 67
$ws.Range("F31").Value = 303
$ws.Range("F36").Value = 938
$ws.Range("F37").Value = 7
$ws.Range("F40").Value = 406
$ws.Range("F43").Value = 2
$ws.Range("F46").Value = 44
$ws.Range("F48").Value = 95

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 31
$ws.Range("G6").Value = 126
$ws.Range("F15").Value = 50
$ws.Range("F19").Value = 19
$ws.Range("F20").Value = 349

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 203

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 203
$ws.Range("F5").Value = 31
$ws.Range("F9").Value = 1151
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202409/rDtevgY01726641860615.jpeg"
$ws.Range("F13").Value = 230
$ws.Range("G15").Value = 126
$ws.Range("F17").Value = 8027
$ws.Range("F18").Value = 9540
$ws.Range("F23").Value = 46
$ws.Range("F25").Value = 1638
$ws.Range("F27").Value = 67
$ws.Range("F28").Value = 303
$ws.Range("F34").Value = 938
$ws.Range("F38").Value = 1396
$ws.Range("F39").Value = 50
$ws.Range("F44").Value = 44
$ws.Range("F47").Value = 19
$ws.Range("F48").Value = 349
$ws.Range("F49").Value = 95
